$wb = $excel.ActiveWorkbook

# Sheet ALC, row 10
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 29999
$ws.Range("J10").Value = 29999
$ws.Range("L10").Value = 29999
$ws.Range("N10").Value = -30585

# Sheet ALC, row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3449.75
$ws.Range("J17").Value = 3403.5898
$ws.Range("L17").Value = 10210.7694
$ws.Range("N17").Value = -10546.7694

# Sheet ALC, row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4024.4092
$ws.Range("I51").Value = 4222.222
$ws.Range("K51").Value = 4222.222
$ws.Range("M51").Value = -3738.222

# Sheet ALC, row 109
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 109684
$ws.Range("J109").Value = 109684
$ws.Range("L109").Value = 109684
$ws.Range("N109").Value = -112458

# Sheet ALC, row 120
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H120").Value = 83836.2
$ws.Range("J120").Value = 83836.2
$ws.Range("L120").Value = 83836.2
$ws.Range("N120").Value = -93512.2

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 6175738.5
$ws.Range("I137").Value = 2178.2
$ws.Range("J137").Value = 13892688
$ws.Range("K137").Value = 6534.599999999999
$ws.Range("L137").Value = 41678064
$ws.Range("M137").Value = -3984.599999999999
$ws.Range("N137").Value = -41683164

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2459.87
$ws.Range("I138").Value = 1187.875
$ws.Range("J138").Value = 2702.1548
$ws.Range("K138").Value = 3563.625
$ws.Range("L138").Value = 8106.464399999999
$ws.Range("M138").Value = 1576.375
$ws.Range("N138").Value = -18386.4644

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1651.5555
$ws.Range("I2").Value = 1271.7778
$ws.Range("K2").Value = 1271.7778
$ws.Range("M2").Value = -1158.7778

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2814.611
$ws.Range("I74").Value = 2803.7058
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 2803.7058
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -1929.7058
$ws.Range("N74").Value = -4748

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2814.611
$ws.Range("I77").Value = 2803.7058
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 14018.529
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -9650.529
$ws.Range("N77").Value = -23736

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1651.5555
$ws.Range("I116").Value = 1271.7778
$ws.Range("K116").Value = 1271.7778
$ws.Range("M116").Value = 1022.2222

# Sheet ARM, row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 49797.5
$ws.Range("J139").Value = 52000
$ws.Range("L139").Value = 52000
$ws.Range("N139").Value = -62280

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1651.5555
$ws.Range("I3").Value = 1271.7778
$ws.Range("K3").Value = 1271.7778
$ws.Range("M3").Value = -1157.7778

# Sheet BSM, row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 447.9091
$ws.Range("I22").Value = 280.77777
$ws.Range("J22").Value = 1200
$ws.Range("K22").Value = 280.77777
$ws.Range("L22").Value = 1200
$ws.Range("M22").Value = -107.77777
$ws.Range("N22").Value = -1546

# Sheet BSM, row 38
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 98830
$ws.Range("J38").Value = 98830
$ws.Range("L38").Value = 98830
$ws.Range("N38").Value = -99662

# Sheet BSM, row 97
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 70476
$ws.Range("I97").Value = 49214
$ws.Range("J97").Value = 113000
$ws.Range("K97").Value = 49214
$ws.Range("L97").Value = 113000
$ws.Range("M97").Value = -48223
$ws.Range("N97").Value = -114982

# Sheet BSM, row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1025.7
$ws.Range("I107").Value = 837.6667
$ws.Range("J107").Value = 1589.8
$ws.Range("K107").Value = 837.6667
$ws.Range("L107").Value = 1589.8
$ws.Range("M107").Value = 1082.3333
$ws.Range("N107").Value = -5429.8

# Sheet CRP, row 20
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 119991
$ws.Range("J20").Value = 119991
$ws.Range("L20").Value = 119991
$ws.Range("N20").Value = -120463

# Sheet CRP, row 30
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H30").Value = 119991
$ws.Range("J30").Value = 119991
$ws.Range("L30").Value = 119991
$ws.Range("N30").Value = -120173

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5918.696
$ws.Range("I31").Value = 2101
$ws.Range("K31").Value = 2101
$ws.Range("M31").Value = -1806

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5918.696
$ws.Range("I34").Value = 2101
$ws.Range("K34").Value = 2101
$ws.Range("M34").Value = -1899

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2777.2559
$ws.Range("J58").Value = 3900
$ws.Range("L58").Value = 3900
$ws.Range("N58").Value = -4306

# Sheet CRP, row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1492.2106
$ws.Range("I107").Value = 736.7778
$ws.Range("J107").Value = 2172.1
$ws.Range("K107").Value = 736.7778
$ws.Range("L107").Value = 2172.1
$ws.Range("M107").Value = 1183.2222
$ws.Range("N107").Value = -6012.1

# Sheet CRP, row 128
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H128").Value = 119991
$ws.Range("J128").Value = 119991
$ws.Range("L128").Value = 119991
$ws.Range("N128").Value = -129951

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2704.7856
$ws.Range("I132").Value = 2612.054
$ws.Range("K132").Value = 7836.162
$ws.Range("M132").Value = -5306.162

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2777.2559
$ws.Range("J136").Value = 3900
$ws.Range("L136").Value = 11700
$ws.Range("N136").Value = -16800

# Sheet CUL, row 70
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2006
$ws.Range("I70").Value = 2006
$ws.Range("K70").Value = 6018
$ws.Range("M70").Value = -5703

# Sheet CUL, row 73
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 2006
$ws.Range("I73").Value = 2006
$ws.Range("K73").Value = 6018
$ws.Range("M73").Value = -4926

# Sheet CUL, row 109
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 3758.2
$ws.Range("I109").Value = 2232.8
$ws.Range("J109").Value = 4266.6665
$ws.Range("K109").Value = 6698.400000000001
$ws.Range("L109").Value = 12799.9995
$ws.Range("M109").Value = -5658.400000000001
$ws.Range("N109").Value = -14879.9995

# Sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3390.9167
$ws.Range("J80").Value = 4499.5
$ws.Range("L80").Value = 4499.5
$ws.Range("N80").Value = -6495.5

# Sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3390.9167
$ws.Range("J83").Value = 4499.5
$ws.Range("L83").Value = 22497.5
$ws.Range("N83").Value = -32481.5

# Sheet GSM, row 100
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 34677.5
$ws.Range("J100").Value = 34677.5
$ws.Range("L100").Value = 34677.5
$ws.Range("N100").Value = -36841.5

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1372.75
$ws.Range("I122").Value = 1372.75
$ws.Range("K122").Value = 4118.25
$ws.Range("M122").Value = -1668.25

# Sheet LTW, row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5694.706
$ws.Range("J7").Value = 5580
$ws.Range("L7").Value = 5580
$ws.Range("N7").Value = -5804

# Sheet LTW, row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4756.125
$ws.Range("I16").Value = 3637.5
$ws.Range("J16").Value = 5874.75
$ws.Range("K16").Value = 3637.5
$ws.Range("L16").Value = 5874.75
$ws.Range("M16").Value = -3467.5
$ws.Range("N16").Value = -6214.75

# Sheet LTW, row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2646.2222
$ws.Range("I22").Value = 1444.6666
$ws.Range("J22").Value = 5049.3335
$ws.Range("K22").Value = 1444.6666
$ws.Range("L22").Value = 5049.3335
$ws.Range("M22").Value = -1149.6666
$ws.Range("N22").Value = -5639.3335

# Sheet LTW, row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2646.2222
$ws.Range("I27").Value = 1444.6666
$ws.Range("J27").Value = 5049.3335
$ws.Range("K27").Value = 1444.6666
$ws.Range("L27").Value = 5049.3335
$ws.Range("M27").Value = -1337.6666
$ws.Range("N27").Value = -5263.3335

# Sheet LTW, row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2712.6
$ws.Range("I93").Value = 2244.4546
$ws.Range("K93").Value = 2244.4546
$ws.Range("M93").Value = -996.4546

# Sheet LTW, row 102
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H102").Value = 115000
$ws.Range("J102").Value = 115000
$ws.Range("L102").Value = 115000
$ws.Range("N102").Value = -121490

# Sheet LTW, row 103
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 109734
$ws.Range("J103").Value = 109734
$ws.Range("L103").Value = 109734
$ws.Range("N103").Value = -112078

# Sheet LTW, row 107
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 7995
$ws.Range("I107").Value = 7995
$ws.Range("K107").Value = 7995
$ws.Range("M107").Value = -6075

# Sheet LTW, row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5694.706
$ws.Range("J126").Value = 5580
$ws.Range("L126").Value = 16740
$ws.Range("N126").Value = -21680

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7602.8667
$ws.Range("J136").Value = 11992.5
$ws.Range("L136").Value = 35977.5
$ws.Range("N136").Value = -41077.5

# Sheet WVR, row 30
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 5958.5
$ws.Range("I30").Value = 5958.5
$ws.Range("K30").Value = 5958.5
$ws.Range("M30").Value = -5851.5

# Sheet WVR, row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 10802.737
$ws.Range("I96").Value = 4867.1665
$ws.Range("J96").Value = 13542.23
$ws.Range("K96").Value = 4867.1665
$ws.Range("L96").Value = 13542.23
$ws.Range("M96").Value = -3494.1665
$ws.Range("N96").Value = -16288.23

# Sheet WVR, row 102
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H102").Value = 85000
$ws.Range("J102").Value = 85000
$ws.Range("L102").Value = 85000
$ws.Range("N102").Value = -91490

# Sheet WVR, row 105
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 30615
$ws.Range("J105").Value = 30615
$ws.Range("L105").Value = 30615
$ws.Range("N105").Value = -37603
